# BalanceTest: relabel the "To/From" column header and switch several
# rows' "To/From" values from location-name strings to numeric balance
# codes (while some stay as (different) location strings).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header relabel: "To/From (Static)" -> "To/From (Location)"
$ws.Range("E1").Value = "To/From (Location)"

# Column E ("To/From") body values
$ws.Range("E3").Value = "Facility 3"
$ws.Range("E4").Value = 101
$ws.Range("E5").Value = "External facility"
$ws.Range("E6").Value = 102
$ws.Range("E7").Value = "Pharmacy"
$ws.Range("E8").Value = "Internal Facility 2"
$ws.Range("E9").Value = 6
$ws.Range("E10").Value = "Final Facility"
$ws.Range("E11").Value = 105
$ws.Range("E12").Value = 103
$ws.Range("E13").Value = 104
$ws.Range("E14").Value = "Imprest"
$ws.Range("E15").Value = "Loc-1"
$ws.Range("E16").Value = "overall ecase"

# Column width / layout adjustments
$ws.Columns.Item(1).ColumnWidth = 4.9166
$ws.Columns.Item(2).ColumnWidth = 10.4166
